# Control the TC Exec via flag
#
# - Row 2 (loginLogoutTest): Execution Flag changes from "yes" to "no"
# - Row 3 (newTest): Execution Flag changes from "no" to "yes", and Priority changes from "2" to "1"
# - Selection / active cell moves from E4 to E8
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle the Execution Flag values
$ws.Range("C2").Value = "no"
$ws.Range("C3").Value = "yes"

# Update Priority for newTest row; keep it stored as text (like the other
# Priority cells) rather than letting Excel auto-convert it to a number.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1"
$ws.Range("D3").NumberFormat = "General"

# Move the active selection
$ws.Range("E8").Select() | Out-Null
